# Mise à jour de l'application
# Append 10 new "Entrainement" rows (J-3 session, date 2025-10-22 / serial 45952)
# to the bottom of the GPS data sheet, mirroring the layout/styling of the
# existing rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 728
$firstNew = 729
$lastNew = 738

# Clone the formatting (date style on col B, centered style on col D, etc.)
# from the last existing row onto the new block before filling in values.
$ws.Range("A728:V728").Copy()
$ws.Range(("A{0}:V{1}" -f $firstNew, $lastNew)).PasteSpecial(-4122)

$rows = @(
    @{ E = "Romain Thunet";     F = "center back";      G = "01:28:40";
       H = 5.94; I = 0.17; J = 5.76; K = 0.15; L = 0.03; M = 0;    N = 0; O = 0;
       P = 3.97; Q = 22.56; R = 4.88; S = 22; T = 4;  U = 18; V = 3 },
    @{ E = "Kamal Bafounta";    F = "center midfield";  G = "01:29:37";
       H = 5.96; I = 0.18; J = 5.77; K = 0.16; L = 0.02; M = 0;    N = 0; O = 0;
       P = 3.94; Q = 22.64; R = 3.97; S = 27; T = 0;  U = 17; V = 1 },
    @{ E = "Naim Dhib";         F = "center midfield";  G = "01:08:16";
       H = 4.51; I = 0.36; J = 4.14; K = 0.25; L = 0.1;  M = 0.02; N = 0; O = 2;
       P = 3.56; Q = 28.99; R = 6.25; S = 62; T = 27; U = 49; V = 25 },
    @{ E = "Yoann Martelat";    F = "center midfield";  G = "01:32:03";
       H = 6.69; I = 0.23; J = 6.46; K = 0.17; L = 0.05; M = 0.01; N = 0; O = 2;
       P = 4.29; Q = 25.67; R = 4.33; S = 17; T = 1;  U = 14; V = 0 },
    @{ E = "Mattheo Haon";      F = "right back";       G = "01:30:50";
       H = 5.86; I = 0.41; J = 5.44; K = 0.36; L = 0.06; M = 0;    N = 0; O = 0;
       P = 3.78; Q = 23.62; R = 4.69; S = 36; T = 6;  U = 31; V = 4 },
    @{ E = "Ilyes Boughanmi";   F = "center forward";   G = "01:28:31";
       H = 5.72; I = 0.24; J = 5.47; K = 0.18; L = 0.07; M = 0;    N = 0; O = 1;
       P = 3.47; Q = 25.7;  R = 5.57; S = 69; T = 23; U = 65; V = 23 },
    @{ E = "Omar Benyounes";    F = "center midfield";  G = "01:30:57";
       H = 6.24; I = 0.31; J = 5.91; K = 0.28; L = 0.04; M = 0;    N = 0; O = 0;
       P = 4.03; Q = 23.87; R = 4.65; S = 38; T = 6;  U = 29; V = 3 },
    @{ E = "Malik Boussaid";    F = "right back";       G = "01:31:30";
       H = 6.35; I = 0.28; J = 6.06; K = 0.25; L = 0.04; M = 0;    N = 0; O = 1;
       P = 3.54; Q = 26.68; R = 5.11; S = 56; T = 8;  U = 52; V = 20 },
    @{ E = "Emmanuel Valey";    F = "left forward";     G = "01:30:18";
       H = 7.79; I = 0.34; J = 7.43; K = 0.3;  L = 0.05; M = 0;    N = 0; O = 0;
       P = 4.11; Q = 24.16; R = 6.32; S = 53; T = 18; U = 47; V = 12 },
    @{ E = "Karahali Souaré";   F = "right forward";    G = "01:29:37";
       H = 6.26; I = 0.33; J = 5.91; K = 0.27; L = 0.08; M = 0;    N = 0; O = 0;
       P = 3.78; Q = 24.42; R = 5.47; S = 72; T = 20; U = 63; V = 16 }
)

$r = $firstNew
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value  = "Entrainement"
    $ws.Cells.Item($r, 2).Value  = 45952
    $ws.Cells.Item($r, 3).Value  = "Global"
    $ws.Cells.Item($r, 4).Value  = "J-3"
    $ws.Cells.Item($r, 5).Value  = $row.E
    $ws.Cells.Item($r, 6).Value  = $row.F
    $ws.Cells.Item($r, 7).Value  = $row.G
    $ws.Cells.Item($r, 8).Value  = $row.H
    $ws.Cells.Item($r, 9).Value  = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $ws.Cells.Item($r, 21).Value = $row.U
    $ws.Cells.Item($r, 22).Value = $row.V
    $r = $r + 1
}

# Match the saved view state (scrolled / selected cell) from the edit.
$excel.ActiveWindow.ScrollRow = 718
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C745").Select()
